$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1698.75
$ws.Range("J19").Value = 1810.75
$ws.Range("L19").Value = 1810.75
$ws.Range("N19").Value = -2160.75
# Row 40
$ws.Range("H40").Value = 2942.7144
$ws.Range("I40").Value = 2870
$ws.Range("J40").Value = 2997.25
$ws.Range("K40").Value = 2870
$ws.Range("L40").Value = 2997.25
$ws.Range("M40").Value = -2695
$ws.Range("N40").Value = -3347.25
# Row 43
$ws.Range("H43").Value = 1626.1666
$ws.Range("J43").Value = 1711.6
$ws.Range("L43").Value = 1711.6
$ws.Range("N43").Value = -1849.6
# Row 70
$ws.Range("H70").Value = 24916.666
$ws.Range("J70").Value = 24916.666
$ws.Range("L70").Value = 74749.99800000001
$ws.Range("N70").Value = -75289.99800000001
# Row 73
$ws.Range("H73").Value = 24916.666
$ws.Range("J73").Value = 24916.666
$ws.Range("L73").Value = 74749.99800000001
$ws.Range("N73").Value = -76621.99800000001
# Row 129
$ws.Range("H129").Value = 946.8333
$ws.Range("I129").Value = 1136.1666
$ws.Range("J129").Value = 883.7222
$ws.Range("K129").Value = 3408.4998
$ws.Range("L129").Value = 2651.1666
$ws.Range("M129").Value = 1591.5002
$ws.Range("N129").Value = -12651.1666
# Row 138
$ws.Range("H138").Value = 3670.65
$ws.Range("I138").Value = 4457.5386
$ws.Range("J138").Value = 2209.2856
$ws.Range("K138").Value = 13372.6158
$ws.Range("L138").Value = 6627.8568
$ws.Range("M138").Value = -8232.6158
$ws.Range("N138").Value = -16907.8568

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 44002.5
$ws.Range("J23").Value = 17999
$ws.Range("L23").Value = 17999
$ws.Range("N23").Value = -18517
# Row 32
$ws.Range("H32").Value = 3354.4822
$ws.Range("I32").Value = 2203.3096
$ws.Range("K32").Value = 2203.3096
$ws.Range("M32").Value = -1916.3096
# Row 45
$ws.Range("H45").Value = 1766.4286
$ws.Range("I45").Value = 1713.4
$ws.Range("K45").Value = 1713.4
$ws.Range("M45").Value = -1336.4
# Row 61
$ws.Range("H61").Value = 3508.7917
$ws.Range("I61").Value = 2711.6843
$ws.Range("J61").Value = 6537.8
$ws.Range("K61").Value = 2711.6843
$ws.Range("L61").Value = 6537.8
$ws.Range("M61").Value = -2499.6843
$ws.Range("N61").Value = -6961.8
# Row 63
$ws.Range("H63").Value = 7495
$ws.Range("J63").Value = 7988.5
$ws.Range("L63").Value = 7988.5
$ws.Range("N63").Value = -9360.5
# Row 66
$ws.Range("H66").Value = 7495
$ws.Range("J66").Value = 7988.5
$ws.Range("L66").Value = 39942.5
$ws.Range("N66").Value = -46806.5
# Row 88
$ws.Range("H88").Value = 2793.0557
$ws.Range("I88").Value = 2136.4443
$ws.Range("K88").Value = 2136.4443
$ws.Range("M88").Value = -1730.4443
# Row 91
$ws.Range("H91").Value = 2793.0557
$ws.Range("I91").Value = 2136.4443
$ws.Range("K91").Value = 2136.4443
$ws.Range("M91").Value = -732.4443000000001
# Row 97
$ws.Range("H97").Value = 620.7143
$ws.Range("J97").Value = 366
$ws.Range("L97").Value = 366
$ws.Range("N97").Value = -1358
# Row 136
$ws.Range("H136").Value = 3508.7917
$ws.Range("I136").Value = 2711.6843
$ws.Range("J136").Value = 6537.8
$ws.Range("K136").Value = 8135.0529
$ws.Range("L136").Value = 19613.4
$ws.Range("M136").Value = -5585.0529
$ws.Range("N136").Value = -24713.4

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1343.2609
$ws.Range("J20").Value = 1017
$ws.Range("L20").Value = 1017
$ws.Range("N20").Value = -1511
# Row 86
$ws.Range("H86").Value = 97027.14
$ws.Range("I86").Value = 1969.7858
$ws.Range("J86").Value = 287141.84
$ws.Range("K86").Value = 1969.7858
$ws.Range("L86").Value = 287141.84
$ws.Range("M86").Value = -846.7858000000001
$ws.Range("N86").Value = -289387.84
# Row 89
$ws.Range("H89").Value = 97027.14
$ws.Range("I89").Value = 1969.7858
$ws.Range("J89").Value = 287141.84
$ws.Range("K89").Value = 9848.929
$ws.Range("L89").Value = 1435709.2
$ws.Range("M89").Value = -4232.929
$ws.Range("N89").Value = -1446941.2
# Row 99
$ws.Range("H99").Value = 1546.6923
$ws.Range("I99").Value = 1512.2222
$ws.Range("J99").Value = 1624.25
$ws.Range("K99").Value = 1512.2222
$ws.Range("L99").Value = 1624.25
$ws.Range("M99").Value = -14.22219999999993
$ws.Range("N99").Value = -4620.25
# Row 105
$ws.Range("H105").Value = 2240.7144
$ws.Range("I105").Value = 2322.875
$ws.Range("J105").Value = 1977.8
$ws.Range("K105").Value = 2322.875
$ws.Range("L105").Value = 1977.8
$ws.Range("M105").Value = -575.875
$ws.Range("N105").Value = -5471.8
# Row 134
$ws.Range("H134").Value = 7807.2144
$ws.Range("I134").Value = 8260.708000000001
$ws.Range("K134").Value = 24782.124
$ws.Range("M134").Value = -22247.124

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 849.6
$ws.Range("I22").Value = 149
$ws.Range("J22").Value = 1316.6666
$ws.Range("K22").Value = 149
$ws.Range("L22").Value = 1316.6666
$ws.Range("M22").Value = 201
$ws.Range("N22").Value = -2016.6666
# Row 70
$ws.Range("H70").Value = 28666.666
$ws.Range("J70").Value = 28666.666
$ws.Range("L70").Value = 28666.666
$ws.Range("N70").Value = -29296.666
# Row 73
$ws.Range("H73").Value = 28666.666
$ws.Range("J73").Value = 28666.666
$ws.Range("L73").Value = 28666.666
$ws.Range("N73").Value = -30850.666
# Row 132
$ws.Range("H132").Value = 2132.3044
$ws.Range("J132").Value = 5001.7144
$ws.Range("L132").Value = 15005.1432
$ws.Range("N132").Value = -20065.1432

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 133
$ws.Range("H133").Value = 15628539
$ws.Range("I133").Value = 50001130
$ws.Range("J133").Value = 4636.273
$ws.Range("K133").Value = 150003390
$ws.Range("L133").Value = 13908.819
$ws.Range("M133").Value = -149998330
$ws.Range("N133").Value = -24028.819

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2450.625
$ws.Range("I80").Value = 1852.5
$ws.Range("J80").Value = 2650
$ws.Range("K80").Value = 1852.5
$ws.Range("L80").Value = 2650
$ws.Range("M80").Value = -854.5
$ws.Range("N80").Value = -4646
# Row 83
$ws.Range("H83").Value = 2450.625
$ws.Range("I83").Value = 1852.5
$ws.Range("J83").Value = 2650
$ws.Range("K83").Value = 9262.5
$ws.Range("L83").Value = 13250
$ws.Range("M83").Value = -4270.5
$ws.Range("N83").Value = -23234
# Row 102
$ws.Range("H102").Value = 2473.4119
$ws.Range("I102").Value = 3039.1428
$ws.Range("K102").Value = 3039.1428
$ws.Range("M102").Value = -1417.1428
# Row 126
$ws.Range("H126").Value = 86668.836
$ws.Range("I126").Value = 3779.111
$ws.Range("K126").Value = 11337.333
$ws.Range("M126").Value = -8867.332999999999
# Row 132
$ws.Range("H132").Value = 3524.889
$ws.Range("I132").Value = 2422.0527
$ws.Range("K132").Value = 7266.158100000001
$ws.Range("M132").Value = -4736.158100000001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1314.6666
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3590
# Row 27
$ws.Range("H27").Value = 1314.6666
$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3214
# Row 46
$ws.Range("H46").Value = 1435.2667
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1435.2667
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1435.2667
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1811.2667
# Row 55
$ws.Range("H55").Value = 237
$ws.Range("I55").Value = 225.8
$ws.Range("J55").Value = 242.09091
$ws.Range("K55").Value = 225.8
$ws.Range("L55").Value = 242.09091
$ws.Range("M55").Value = -52.80000000000001
$ws.Range("N55").Value = -588.09091
# Row 68
$ws.Range("H68").Value = 3097.75
$ws.Range("I68").Value = 2826
$ws.Range("K68").Value = 2826
$ws.Range("M68").Value = -2077
# Row 71
$ws.Range("H71").Value = 3097.75
$ws.Range("I71").Value = 2826
$ws.Range("K71").Value = 14130
$ws.Range("M71").Value = -10386
# Row 74
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996
# Row 77
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984
# Row 82
$ws.Range("H82").Value = 3728
$ws.Range("I82").Value = 1875
$ws.Range("J82").Value = 4191.25
$ws.Range("K82").Value = 1875
$ws.Range("L82").Value = 4191.25
$ws.Range("M82").Value = -1514
$ws.Range("N82").Value = -4913.25
# Row 85
$ws.Range("H85").Value = 3728
$ws.Range("I85").Value = 1875
$ws.Range("J85").Value = 4191.25
$ws.Range("K85").Value = 1875
$ws.Range("L85").Value = 4191.25
$ws.Range("M85").Value = -627
$ws.Range("N85").Value = -6687.25
# Row 136
$ws.Range("H136").Value = 3440.2334
$ws.Range("I136").Value = 2357.6191
$ws.Range("K136").Value = 7072.8573
$ws.Range("M136").Value = -4522.8573

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 28000
$ws.Range("J64").Value = 28000
$ws.Range("L64").Value = 28000
$ws.Range("N64").Value = -28496
# Row 67
$ws.Range("H67").Value = 28000
$ws.Range("J67").Value = 28000
$ws.Range("L67").Value = 28000
$ws.Range("N67").Value = -29716
